$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.839.64"
$ws.Range("E2").Value = "  -0.74%  "

$ws.Range("D3").Value = "3.452.49"

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'574.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.20%  "

$ws.Range("D6").Value = "'159.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.20%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "3.450.58"
$ws.Range("E8").Value = "  -1.31%  "

$ws.Range("D9").Value = "'0.575"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.17%  "

$ws.Range("D10").Value = "'7.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.15%  "

$ws.Range("D11").Value = "'0.121"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.39%  "

$ws.Range("D12").Value = "'0.439"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.08%  "

$ws.Range("D13").Value = "4.046.80"
$ws.Range("E13").Value = "  -1.17%  "

$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("D15").Value = "'27.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.26%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "64.894.74"
$ws.Range("E16").Value = "  -0.66%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000173"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -10.91%  "

$ws.Range("D18").Value = "3.459.09"
$ws.Range("E18").Value = "  -0.91%  "

$ws.Range("D19").Value = "'6.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.30%  "

$ws.Range("D20").Value = "'13.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.08%  "

$ws.Range("D21").Value = "'376.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.44%  "

$ws.Range("D22").Value = "'7.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.71%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").Value = "'72.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.67%  "

$ws.Range("E25").Value = "  -3.55%  "

$ws.Range("D26").Value = "'0.0000120"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.62%  "

$ws.Range("D27").Value = "'9.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.06%  "

$ws.Range("E28").Value = "  -0.24%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("D30").Value = "'1.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.55%  "

$ws.Range("E31").Value = "  -2.08%  "

$ws.Range("E32").Value = "  -2.45%  "

$ws.Range("D33").Value = "'23.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.68%  "

$ws.Range("D34").Value = "'6.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.14%  "

$ws.Range("E35").Value = "  -4.46%  "

$ws.Range("D36").Value = "'161.22"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.82%  "

$ws.Range("E37").Value = "  -2.98%  "

$ws.Range("D38").Value = "2.898.27"
$ws.Range("E38").Value = "  -4.11%  "

$ws.Range("D39").Value = "'0.0748"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.32%  "

$ws.Range("D40").Value = "'26.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.40%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'43.02"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'4.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.92%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.787"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.47%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'26.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.30%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'6.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.62%  "

$ws.Range("E46").Value = "  -3.73%  "

$ws.Range("E47").Value = "  +8.04%  "

$ws.Range("D48").Value = "'321.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.37%  "

$ws.Range("E49").Value = "  -3.56%  "

$ws.Range("D50").Value = "'6.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.48%  "

$ws.Range("D51").Value = "'0.844"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.56%  "
